$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "victory"
$ws.Range("B5").Value = "VICTORY"

$ws.Range("B5").Select()
